$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")

$newRow = 81

# Column A holds plain text dates ("YYYY-MM-DD"), stored as shared strings
# with no special cell style/number-format in the source workbook. Writing
# a date-shaped literal straight into Value/Formula makes Excel's
# autodetection coerce it into a real date (and silently stamps a new
# date-flavoured cell style onto the cell). To land an actual text value
# that round-trips as a shared string with the default style - matching
# the original file's convention - compute the text with TEXT() and copy
# it back in as a value via PasteSpecial (values only), which bypasses
# the literal-entry autodetection path entirely.
$ws.Cells.Item($newRow, 1).Formula = '=TEXT(44062,"yyyy-mm-dd")'
$ws.Cells.Item($newRow, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4163) | Out-Null

$ws.Cells.Item($newRow, 2).Value = 537031
$ws.Cells.Item($newRow, 3).Value = 591637
$ws.Cells.Item($newRow, 4).Value = 82884
$ws.Cells.Item($newRow, 5).Value = 58481
$ws.Cells.Item($newRow, 6).Value = 26.09
